$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.186.86"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.522.77"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.52"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.43"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.523.39"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.12"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.121.98"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.63"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.518.44"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.189.06"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("E19").Value = "  -3.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.03"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.63"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "385.40"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.665.69"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.02"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000116"
$ws.Range("E27").Value = "  +2.31%  "
$ws.Range("E28").Value = "  -2.51%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.50"
$ws.Range("E30").Value = "  -3.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.43"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("E32").Value = "  -2.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.533.47"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.58"
$ws.Range("E35").Value = "  -2.38%  "
$ws.Range("E36").Value = "  +1.11%  "
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.92"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "161.04"
$ws.Range("E40").Value = "  -4.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0788"
$ws.Range("E41").Value = "  -2.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.814"
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.63"
$ws.Range("E43").Value = "  +3.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.64"
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("E46").Value = "  -4.87%  "
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.62"
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.478.01"
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.907"
$ws.Range("E51").Value = "  -0.55%  "
